$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update estoque_atualizado (column G) values per the diff
$ws.Range("G2").Value = -455
$ws.Range("G3").Value = -455
$ws.Range("G5").Value = -455
$ws.Range("G10").Value = 2
